$wb = $excel.ActiveWorkbook

# --- Overview sheet: status message changes from "in sync" to "not in sync" ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: not in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: not in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: not in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: not in sync with en-US"

# Columns E/F widened (report regenerated with longer status text)
$wsOverview.Range("E1").ColumnWidth = 32.65
$wsOverview.Range("F1").ColumnWidth = 32.65

# --- zh-cn sheet: status message + new handback datetime for the first file ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Handed back: not in sync with en-US"
$wsZhCn.Range("C3").Value = "Handed back: not in sync with en-US"
$wsZhCn.Range("K2").Value = "2016-09-07 14:45:03"

# Status column widened to match
$wsZhCn.Range("C1").ColumnWidth = 32.65

# --- de-de sheet: status message + new handback datetime for the first file ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Handed back: not in sync with en-US"
$wsDeDe.Range("C3").Value = "Handed back: not in sync with en-US"
$wsDeDe.Range("K2").Value = "2016-09-07 14:45:52"

# Status column widened to match
$wsDeDe.Range("C1").ColumnWidth = 32.65
